$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) text to match the corrected MOH template ---
$ws.Range("A1").Value = "שם הספק"
$ws.Range("B1").Value = "ח""פ ספק "
$ws.Range("C1").Value = "מספר משרד הבריאות"
$ws.Range("D1").Value = "תאריך"
$ws.Range("E1").Value = "מס.רכב"
$ws.Range("F1").Value = "שם הנהג"
$ws.Range("G1").Value = "טלפון נהג"
$ws.Range("H1").Value = "לקוח"
$ws.Range("I1").Value = "סוג לקוח (קמעונאי,מפעל/מחסן)"
$ws.Range("J1").Value = "קוד עיר"
$ws.Range("K1").Value = "כתובת"
$ws.Range("L1").Value = "ח""פ לקוח `nאו מספר אישור משרד הבריאות במקרים בהם המשלוח הוא למפעל מאושר"
$ws.Range("M1").Value = "מספר סניף הרשת"
$ws.Range("N1").Value = "מספר תעודת משלוח"
$ws.Range("O1").Value = "בשר בהמות גולמי"
$ws.Range("P1").Value = "בשר בהמות מיבוא קפוא"
$ws.Range("Q1").Value = "בשר בהמות מעובד"
$ws.Range("R1").Value = "עוף גולמי (עוף שחוט)"
$ws.Range("S1").Value = "עוף מעובד"
$ws.Range("T1").Value = "דגים גולמי (מקומי)"
$ws.Range("U1").Value = "דגים יבוא"
$ws.Range("V1").Value = "דגים מעובדים"
$ws.Range("W1").Value = "מוצרים מוכנים לאכילה"
$ws.Range("X1").Value = "נוסף א"
$ws.Range("Y1").Value = "נוסף ב"
$ws.Range("Z1").Value = "סה""כ קרטונים"
$ws.Range("AA1").Value = "סה""כ משקל"
$ws.Range("AB1").Value = "סבב יומי"
$ws.Range("AC1").Value = "קוד ביטול דיווח משלוח`n(למקרים בהם נדרש לבטל תעודת משלוח שדווחה ולא יצאה מהמפעל לשיווק"
$ws.Range("AD1").Value = "משווק באמצעות"

# --- Update delivery-timestamp column (D) for rows 2-14 to corrected values ---
$ws.Range("D2").Value = 45944.8221248004
$ws.Range("D3").Value = 45944.82212480129
$ws.Range("D4").Value = 45944.822124801634
$ws.Range("D5").Value = 45944.822124801954
$ws.Range("D6").Value = 45944.82212480222
$ws.Range("D7").Value = 45944.82212480263
$ws.Range("D8").Value = 45944.822124802886
$ws.Range("D9").Value = 45944.82212480316
$ws.Range("D10").Value = 45944.82212480336
$ws.Range("D11").Value = 45944.822124803606
$ws.Range("D12").Value = 45944.8221248038
$ws.Range("D13").Value = 45944.822124804
$ws.Range("D14").Value = 45944.822124804195

# --- Remove stray data that was in the wrong columns (M = col 13, W = col 23) ---
$ws.Range("M2:M14").ClearContents()
$ws.Range("W2:W14").ClearContents()

Write-Host "Template headers corrected; column mapping fixed (M/W cleared); D timestamps updated."
